$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the new row's cells as text so that date/time-like
# strings (e.g. "2025-03-11", "22:24:10") are stored as literal text
# instead of being auto-converted into date/time serial numbers.
$ws.Range("A3:D3").NumberFormat = "@"

$ws.Range("A3").Value = "2025-03-11"
$ws.Range("B3").Value = "Vasanth Kumar"
$ws.Range("C3").Value = "22:24:10"
$ws.Range("D3").Value = "22:24:19"

# Restore the default (unstyled) cell style so the new row matches the
# plain, style-less formatting of the existing data rows.
$ws.Range("A3:D3").Style = "Normal"
